$d = $word.ActiveDocument

# --- "Estudiante 1 Cod XXXX" -> "Andres Mendoza 202012676" ---
$d.Content.Find.Execute("Estudiante 1 Cod XXXX", $false, $false, $false, $false, $false, $true, 1, $false, "Andres Mendoza 202012676", 2)

# --- "Estudiante 2 Cod XXXX" -> "Daniela Alvarez 202020209" ---
$d.Content.Find.Execute("Estudiante 2 Cod XXXX", $false, $false, $false, $false, $false, $true, 1, $false, "Daniela Alvarez 202020209", 2)

# --- Remove the now-empty paragraph that used to sit below the two student lines ---
$p = $d.Paragraphs(4).Range
$p.Delete()
